# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped figures (GitHub Actions data refresh).
#
# All values in these two columns are stored as plain text in the workbook
# (e.g. "51.852.95", "  -0.55%  "). Several of the new Price figures look like
# ordinary decimal numbers (e.g. "352.35"); assigning such a string straight to
# .Value would let Excel auto-convert it to a numeric cell, which would both
# change the cells stored type and round-trip the value through binary
# floating point (e.g. 352.35 -> 352.35000000000002). To keep those cells as
# text - exactly like the rest of the column - a leading apostrophe is used,
# mirroring how a user forces "text that looks like a number" into a cell in
# the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.852.95'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.954.52'
$ws.Range("E3").Value = '  +2.74%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''352.35'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '''111.88'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("D10").Value = '''39.64'
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").Value = '''0.0902'
$ws.Range("E11").Value = '  +6.14%  '
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").Value = '''19.86'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").Value = '''8.09'
$ws.Range("E14").Value = '  +2.60%  '
$ws.Range("D15").Value = '3.421.42'
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("D16").Value = '2.955.58'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '51.980.50'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").Value = '''7.76'
$ws.Range("E19").Value = '  +1.39%  '
$ws.Range("D20").Value = '''14.49'
$ws.Range("E20").Value = '  +6.51%  '
$ws.Range("E21").Value = '  -2.58%  '
$ws.Range("D22").Value = '0.0₃0992'
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("D23").Value = '''71.54'
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("D24").Value = '''273.22'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").Value = '''0.182'
$ws.Range("E26").Value = '  +10.86%  '
$ws.Range("D27").Value = '''27.45'
$ws.Range("E27").Value = '  +3.31%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("D29").Value = '''7.46'
$ws.Range("E29").Value = '  +18.91%  '
$ws.Range("D30").Value = '''0.111'
$ws.Range("E30").Value = '  +23.89%  '
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("D32").Value = '''6.38'
$ws.Range("E32").Value = '  +8.73%  '
$ws.Range("D33").Value = '''37.79'
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("D34").Value = '''53.03'
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").Value = '''0.0450'
$ws.Range("E35").Value = '  -1.23%  '
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").Value = '''3.40'
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").Value = '''1.84'
$ws.Range("E38").Value = '  -12.05%  '
$ws.Range("D39").Value = '''18.91'
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("E41").Value = '  +1.16%  '
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("D43").Value = '''23.77'
$ws.Range("E43").Value = '  +5.25%  '
$ws.Range("E44").Value = '  -1.81%  '
$ws.Range("D45").Value = '''3.56'
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").Value = '2.167.12'
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").Value = '''114.25'
$ws.Range("E48").Value = '  -6.86%  '
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("E50").Value = '  +5.67%  '
$ws.Range("D51").Value = '''0.932'
$ws.Range("E51").Value = '  -3.53%  '
